# customers.xlsx: "Update points 09876543 -> 120.00"
#
# The existing row 61 held the special customer phone "09876543" (stored as
# text, with its leading zero) together with its current point total.
# The update re-records the running point total for that same customer by
# converting the old row's phone into a plain number (losing the leading
# zero, as happens when the value is re-entered as numeric) and appending a
# brand-new row 62 that keeps the correctly formatted text phone number
# "09876543" together with the updated point total of 120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61: phone becomes a plain number (9876543) - total_points (C61) stays 0.
$ws.Range("A61").Value = 9876543

# Row 62 (new): text phone "09876543" (leading apostrophe forces text so the
# leading zero survives), blank birthday, and the updated point total.
$ws.Range("A62").Formula = "'09876543"
$ws.Range("B62").NumberFormat = "@"
$ws.Range("C62").Value = 120
